$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.10122516985340724
$ws.Range("A2").Value = -0.0099999994882757903
$ws.Range("A3").Value = -0.0089999994910421321
$ws.Range("A4").Value = 0.28399404951267115
$ws.Range("A5").Value = -0.0059999994994042183
$ws.Range("A6").Value = -0.03851028622023378
$ws.Range("A7").Value = -0.019999999402301682
$ws.Range("A8").Value = -0.019999999396460133
$ws.Range("A9").Value = -0.0059999994632420339
$ws.Range("A10").Value = -0.0059999994554829073
$ws.Range("A11").Value = -0.0044999994625847251
$ws.Range("A12").Value = -0.0059999994529409406
$ws.Range("A13").Value = -0.0059999994449873029
$ws.Range("A14").Value = -0.01199999941051022
$ws.Range("A15").Value = -0.0059999994410606661
$ws.Range("A16").Value = -0.005999999439517012
$ws.Range("A17").Value = -0.0059999994374271282
$ws.Range("A18").Value = -0.008999999420990612
$ws.Range("A19").Value = -0.0089999994933167571
$ws.Range("A20").Value = -0.0089999994892711044
$ws.Range("A21").Value = -0.0089999994887470791
$ws.Range("A22").Value = -0.0089999994884486512
$ws.Range("A23").Value = -0.0089999994873046774
$ws.Range("A24").Value = -0.041999999302794855
$ws.Range("A25").Value = -0.041999999299422441
$ws.Range("A26").Value = -0.0059999994783517252
$ws.Range("A27").Value = -0.0059999994766650744
$ws.Range("A28").Value = -0.0059999994688757496
$ws.Range("A29").Value = -0.011999999431608899
$ws.Range("A30").Value = -0.0092085901982557061
$ws.Range("A31").Value = -0.01074312532871069
$ws.Range("A32").Value = -0.020999999376665635
$ws.Range("A33").Value = -0.0059999994570194559
